# Add season record (Wins / Losses / Ties) columns to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in columns AD, AE, AF (row 1). Copy the formatting of the
# existing header style (bold / bordered / centered) from the last header
# cell (AC1) onto the new header cells before setting their values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team's 1999 season record applied to every data row (rows 2 through 41).
$wins = 98
$losses = 64
$ties = 0

for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
